$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "No aplicable"
$ws.Range("A2").Value = "Equipamientos de salud(ambulatorio, centro de salud, hospital,...)"
$ws.Range("A3").Value = "Local comercial"
$ws.Range("A4").Value = "Oficinas(incluye también el resto de servicios)"
$ws.Range("A5").Value = "Equipamientos de bienestar social(club de ancianos, centro de servicios sociales, centro de día,...)"
$ws.Range("A6").Value = "Local agrario"
$ws.Range("A7").Value = "Equipamientos educativos(colegio, facultad, guardería, escuela,...)"
$ws.Range("A8").Value = "Local industrial"
$ws.Range("A9").Value = "Equipamientos culturales o deportivos(teatro, cine, museo, sala de exposiciones, polideportivo,...)"

$wb.Save()
